$d = $word.ActiveDocument

# Change 1: remove the tracked-insertion leading space at the very start of
# the document (the <w:ins> wrapped run). TrackRevisions is off, so deleting
# this single character removes it outright rather than recording a w:del.
$lead = $d.Range(0, 1)
$lead.Delete()

# Change 2: " and cleaned, the sooner " -> " cleaned, the sooner "
$d.Content.Find.Execute(" and cleaned, the sooner ", $false, $false, $false, $false, $false, $true, 1, $false, " cleaned, the sooner ", 2)

# Change 3: "so that , if there" -> "so that, if there"
$d.Content.Find.Execute("so that , if there", $false, $false, $false, $false, $false, $true, 1, $false, "so that, if there", 2)

# Change 4: "to leave the mold to cool a little outside, one well reheated" ->
# "to leave the mold die down a little outside, once well reheated"
$d.Content.Find.Execute("to leave the mold to cool a little outside, one well reheated", $false, $false, $false, $false, $false, $true, 1, $false, "to leave the mold die down a little outside, once well reheated", 2)
